# Applies the Vietnamese-translation edits described by the diff to the
# "[TEMPLATE] Affiliate email - invite to seminar" document's English
# source section (adding VI text where none existed) and to two
# already-translated Vietnamese sentences further down in the document.

$d = $word.ActiveDocument

# wdReplaceOne = 1, wdReplaceAll = 2, wdFindStop = 0

# --- 1) English section headline / greeting / intro paragraph ----------
$d.Content.Find.Execute(
    "You’re invited to our Deriv Partner Seminar", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Lời mời tham dự Hội thảo Đối tác của Deriv", 2) | Out-Null

$d.Content.Find.Execute(
    "Dear [PARTNER NAME], ", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Xin chào [PARTNER NAME], ", 2) | Out-Null

$d.Content.Find.Execute(
    "We’re excited to let you know that the Deriv Affiliate team will be in [CITY] in [MONTH] to meet with you, our valued partners!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Chúng tôi rất vui mừng thông báo đến bạn đội ngũ Tiếp thị liên kết của Deriv sẽ có mặt tại thành phố [CITY] vào tháng [MONTH] để gặp bạn - đối tác quan trọng của chúng tôi!",
    2) | Out-Null

# --- 2) English section: one-day seminar paragraph ----------------------
$d.Content.Find.Execute(
    "In this one-day seminar, we’ll be providing technical and marketing support, offering the opportunity to network with other partners over a delicious lunch as well as listening to your feedback about our partnership programmes. This is your chance to get your voice heard, which will help us plan future efforts to support you better. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Trong suốt một ngày diễn ra hội thảo, chúng tôi sẽ mang đến bạn nhiều nội dung hỗ trợ về mặt kỹ thuật và marketing, tạo cơ hội kết nối giữa bạn với các đối tác khác qua bữa ăn trưa ngon miệng, đồng thời lắng nghe phản hồi của bạn về các chương trình đối tác của chúng tôi. Đây là cơ hội để bạn có thể chia sẻ những ý kiến đóng góp của mình và điều này sẽ giúp chúng tôi lên kế hoạch để có thể hỗ trợ bạn tốt hơn nữa trong thời gian tới. ",
    2) | Out-Null

# --- 3) English section: RSVP line tail ---------------------------------
$d.Content.Find.Execute(
    ". Please note that attendance is confirmed on a first come, first served basis. We look forward to seeing you there!",
    $true, $false, $false, $false, $false, $true, 1, $false,
    ". Please note that attendance is confirmed on a first come, first served basis. Chúng tôi rất mong được gặp bạn tại sự kiện!",
    2) | Out-Null

# --- 4) "Send my details" button (both occurrences in the document) -----
$d.Content.Find.Execute(
    "Send my details", $true, $false, $false, $false, $false, $true, 1,
    $false, "Gửi thông tin của tôi", 2) | Out-Null

# --- 5) English section: "contact us via" + " or " between hyperlinks ---
$d.Content.Find.Execute(
    "If you have any questions, please contact us via ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua ", 2) | Out-Null

# The standalone " or " between the "live chat" and "WhatsApp" hyperlinks is
# the first of several " or " substrings in the document (a second, unrelated
# one follows later in the same paragraph, and a third lives in the French
# section), so only replace the first hit. A plain Find/Replace right after a
# hyperlink run picks up that hyperlink's character formatting (color +
# underline) for the replacement text, so instead: locate the run with Find
# (no replace), insert the Vietnamese text immediately before the following
# hyperlink (inheriting the plain/neighbouring run's formatting instead) and
# then delete the original " or " span.
$r = $d.Content
$r.Find.Execute(" or ", $true, $false, $false, $false, $false, $true, 1,
    $false, "", 0) | Out-Null
$orStart = $r.Start
$orEnd = $r.End
$insertPoint = $d.Range($orEnd, $orEnd)
$insertPoint.InsertBefore(" hoặc ")
$d.Range($orStart, $orEnd).Delete() | Out-Null

# --- 6) Existing Vietnamese section: tighten two sentences --------------
# Match the FULL original paragraph (not just the tail clause) so this only
# hits the pre-existing Vietnamese-section sentence (which says "... ket noi
# CHO ban ...") and not the newly Vietnamese-ized English-section sentence
# from step 2 above (which now says "... ket noi GIUA ban ..." and must keep
# its own, longer, ending untouched).
$d.Content.Find.Execute(
    "Trong suốt một ngày diễn ra hội thảo, chúng tôi sẽ mang đến bạn nhiều nội dung hỗ trợ về mặt kỹ thuật và marketing, tạo cơ hội kết nối cho bạn với các đối tác khác qua bữa ăn trưa ngon miệng, đồng thời lắng nghe phản hồi của bạn về các chương trình đối tác của chúng tôi. Đây là cơ hội để bạn có thể chia sẻ những ý kiến đóng góp của mình và điều này sẽ giúp chúng tôi lên kế hoạch để có thể hỗ trợ bạn tốt hơn nữa trong thời gian tới. ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Trong suốt một ngày diễn ra hội thảo, chúng tôi sẽ mang đến bạn nhiều nội dung hỗ trợ về mặt kỹ thuật và marketing, tạo cơ hội kết nối cho bạn với các đối tác khác qua bữa ăn trưa ngon miệng, đồng thời lắng nghe phản hồi của bạn về các chương trình đối tác của chúng tôi. Đây là cơ hội để bạn có thể chia sẻ ý kiến của mình đồng thời cũng sẽ giúp chúng tôi lên kế hoạch để có thể hỗ trợ bạn tốt hơn nữa trong thời gian tới. ",
    2) | Out-Null

$d.Content.Find.Execute(
    "Vui lòng phản hồi chúng tôi bằng cách gửi đơn đăng ký trước ngày ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Vui lòng xác nhận tham gia bằng cách gửi đơn đăng ký trước ngày ", 2) |
    Out-Null

$d.Content.Find.Execute(
    "Gửi thông tin cá nhân", $true, $false, $false, $false, $false, $true,
    1, $false, "Gửi thông tin của tôi", 2) | Out-Null
